$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet lists, per language, 4 file-format rows (txt/xml/json/html).
# Add a new language "spa" following the same 4-row block pattern used
# for every other language (see rows 22:25 for "tam").
$srcRange = $ws.Range("A22:D25")
$destRange = $ws.Range("A26:D29")
$srcRange.Copy($destRange)

$ws.Range("A26").Value = "spa"
$ws.Range("A27").Value = "spa"
$ws.Range("A28").Value = "spa"
$ws.Range("A29").Value = "spa"

$ws.Range("A26:D29").Select()
